# Auto-generated Excel COM-interop edit script
# Applies numeric corrections to the Leve profit-tracking columns
# (currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ /
#  LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ) across the
# ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR sheets, matching the refreshed
# market-board snapshot from the scheduled Odin_Profits runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 1550.5
$ws.Range("I11").Value = 1550.5
$ws.Range("K11").Value = 1550.5
$ws.Range("M11").Value = -1410.5
$ws.Range("H53").Value = 474.33334
$ws.Range("J53").Value = 483
$ws.Range("L53").Value = 483
$ws.Range("N53").Value = -1757
$ws.Range("H140").Value = 197887.33
$ws.Range("J140").Value = 231569.42
$ws.Range("L140").Value = 231569.42
$ws.Range("N140").Value = -241929.42
$ws.Range("H141").Value = 3637.6
$ws.Range("I141").Value = 3662.6667
$ws.Range("K141").Value = 10988.0001
$ws.Range("M141").Value = -5808.000100000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4559.9
$ws.Range("I2").Value = 2137.375
$ws.Range("K2").Value = 2137.375
$ws.Range("M2").Value = -2024.375
$ws.Range("H5").Value = 370.42856
$ws.Range("I5").Value = 370.42856
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 370.42856
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -258.42856
$ws.Range("N5").ClearContents() | Out-Null
$ws.Range("H32").Value = 7026265
$ws.Range("I32").Value = 6173999
$ws.Range("K32").Value = 6173999
$ws.Range("M32").Value = -6173712
$ws.Range("H61").Value = 5160.0547
$ws.Range("I61").Value = 6949.9165
$ws.Range("K61").Value = 6949.9165
$ws.Range("M61").Value = -6737.9165
$ws.Range("H102").Value = 3089.5715
$ws.Range("I102").Value = 3058
$ws.Range("K102").Value = 3058
$ws.Range("M102").Value = -1436
$ws.Range("H116").Value = 4559.9
$ws.Range("I116").Value = 2137.375
$ws.Range("K116").Value = 2137.375
$ws.Range("M116").Value = 156.625
$ws.Range("H132").Value = 545985
$ws.Range("I132").Value = 634440.4399999999
$ws.Range("J132").Value = 94862.5
$ws.Range("K132").Value = 1903321.32
$ws.Range("L132").Value = 284587.5
$ws.Range("M132").Value = -1900791.32
$ws.Range("N132").Value = -289647.5
$ws.Range("H136").Value = 5160.0547
$ws.Range("I136").Value = 6949.9165
$ws.Range("K136").Value = 20849.7495
$ws.Range("M136").Value = -18299.7495

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4559.9
$ws.Range("I3").Value = 2137.375
$ws.Range("K3").Value = 2137.375
$ws.Range("M3").Value = -2023.375
$ws.Range("H4").Value = 370.42856
$ws.Range("I4").Value = 370.42856
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 370.42856
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -255.42856
$ws.Range("N4").ClearContents() | Out-Null
$ws.Range("H94").Value = 4454.5
$ws.Range("I94").Value = 1708.6
$ws.Range("K94").Value = 1708.6
$ws.Range("M94").Value = -1257.6
$ws.Range("H99").Value = 6678.387
$ws.Range("I99").Value = 3467.4
$ws.Range("K99").Value = 3467.4
$ws.Range("M99").Value = -1969.4
$ws.Range("H105").Value = 3708.4167
$ws.Range("I105").Value = 3687.75
$ws.Range("J105").Value = 3749.75
$ws.Range("K105").Value = 3687.75
$ws.Range("L105").Value = 3749.75
$ws.Range("M105").Value = -1940.75
$ws.Range("N105").Value = -7243.75
$ws.Range("H134").Value = 2008550.1
$ws.Range("I134").Value = 2786239.5
$ws.Range("K134").Value = 8358718.5
$ws.Range("M134").Value = -8356183.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 35719572
$ws.Range("I16").Value = 62504184
$ws.Range("J16").Value = 6750
$ws.Range("K16").Value = 62504184
$ws.Range("L16").Value = 6750
$ws.Range("M16").Value = -62503897
$ws.Range("N16").Value = -7324
$ws.Range("H38").Value = 6347.25
$ws.Range("I38").Value = 6799.6665
$ws.Range("J38").Value = 4990
$ws.Range("K38").Value = 6799.6665
$ws.Range("L38").Value = 4990
$ws.Range("M38").Value = -6422.6665
$ws.Range("N38").Value = -5744
$ws.Range("H46").Value = 6347.25
$ws.Range("I46").Value = 6799.6665
$ws.Range("J46").Value = 4990
$ws.Range("K46").Value = 6799.6665
$ws.Range("L46").Value = 4990
$ws.Range("M46").Value = -6588.6665
$ws.Range("N46").Value = -5412
$ws.Range("H107").Value = 2226.7222
$ws.Range("I107").Value = 1649.5385
$ws.Range("K107").Value = 1649.5385
$ws.Range("M107").Value = 270.4614999999999
$ws.Range("H113").Value = 35719572
$ws.Range("I113").Value = 62504184
$ws.Range("J113").Value = 6750
$ws.Range("K113").Value = 62504184
$ws.Range("L113").Value = 6750
$ws.Range("M113").Value = -62502014
$ws.Range("N113").Value = -11090
$ws.Range("H132").Value = 12459
$ws.Range("I132").Value = 8577.125
$ws.Range("J132").Value = 20222.75
$ws.Range("K132").Value = 25731.375
$ws.Range("L132").Value = 60668.25
$ws.Range("M132").Value = -23201.375
$ws.Range("N132").Value = -65728.25
$ws.Range("H134").Value = 58829580
$ws.Range("I134").Value = 76927656
$ws.Range("K134").Value = 230782968
$ws.Range("M134").Value = -230780433

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H51").Value = 1639.2
$ws.Range("I51").Value = 799.25
$ws.Range("K51").Value = 2397.75
$ws.Range("M51").Value = -1937.75
$ws.Range("H132").Value = 3326.889
$ws.Range("I132").Value = 988
$ws.Range("J132").Value = 4496.3335
$ws.Range("K132").Value = 8892
$ws.Range("L132").Value = 40467.0015
$ws.Range("M132").Value = -6362
$ws.Range("N132").Value = -45527.0015
$ws.Range("H133").Value = 10555
$ws.Range("J133").Value = 10555
$ws.Range("L133").Value = 31665
$ws.Range("N133").Value = -41785
$ws.Range("H134").Value = 43495908
$ws.Range("I134").Value = 47635400
$ws.Range("K134").Value = 142906200
$ws.Range("M134").Value = -142901130

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 808053.5600000001
$ws.Range("I102").Value = 1187307.2
$ws.Range("J102").Value = 7406.9443
$ws.Range("K102").Value = 1187307.2
$ws.Range("L102").Value = 7406.9443
$ws.Range("M102").Value = -1185685.2
$ws.Range("N102").Value = -10650.9443

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 18421.285
$ws.Range("I7").Value = 16533.166
$ws.Range("J7").Value = 29750
$ws.Range("K7").Value = 16533.166
$ws.Range("L7").Value = 29750
$ws.Range("M7").Value = -16421.166
$ws.Range("N7").Value = -29974
$ws.Range("H16").Value = 2262.1304
$ws.Range("I16").Value = 788.5238000000001
$ws.Range("K16").Value = 788.5238000000001
$ws.Range("M16").Value = -618.5238000000001
$ws.Range("H36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("N36").ClearContents() | Out-Null
$ws.Range("H40").Value = 6195.6206
$ws.Range("I40").Value = 6166.2607
$ws.Range("K40").Value = 6166.2607
$ws.Range("M40").Value = -6030.2607
$ws.Range("H93").Value = 2534.8235
$ws.Range("I93").Value = 2572.8667
$ws.Range("J93").Value = 2249.5
$ws.Range("K93").Value = 2572.8667
$ws.Range("L93").Value = 2249.5
$ws.Range("M93").Value = -1324.8667
$ws.Range("N93").Value = -4745.5
$ws.Range("H122").Value = 4397.577
$ws.Range("I122").Value = 3442.5789
$ws.Range("K122").Value = 10327.7367
$ws.Range("M122").Value = -7877.736699999999
$ws.Range("H126").Value = 18421.285
$ws.Range("I126").Value = 16533.166
$ws.Range("J126").Value = 29750
$ws.Range("K126").Value = 49599.49800000001
$ws.Range("L126").Value = 89250
$ws.Range("M126").Value = -47129.49800000001
$ws.Range("N126").Value = -94190
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents() | Out-Null
$ws.Range("H132").Value = 6678.2
$ws.Range("I132").Value = 5916.6
$ws.Range("K132").Value = 17749.8
$ws.Range("M132").Value = -15219.8
$ws.Range("H136").Value = 38471404
$ws.Range("I136").Value = 125015020
$ws.Range("K136").Value = 375045060
$ws.Range("M136").Value = -375042510

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 3175.25
$ws.Range("I96").Value = 3275.25
$ws.Range("J96").Value = 3075.25
$ws.Range("K96").Value = 3275.25
$ws.Range("L96").Value = 3075.25
$ws.Range("M96").Value = -1902.25
$ws.Range("N96").Value = -5821.25
$ws.Range("H100").Value = 1442.7142
$ws.Range("I100").Value = 1260
$ws.Range("K100").Value = 2520
$ws.Range("M100").Value = -1979
$ws.Range("H107").Value = 2013.2593
$ws.Range("I107").Value = 1680.4286
$ws.Range("K107").Value = 5041.2858
$ws.Range("M107").Value = -3121.2858
$ws.Range("H122").Value = 7539.1353
$ws.Range("I122").Value = 3678.7144
$ws.Range("K122").Value = 11036.1432
$ws.Range("M122").Value = -8586.143199999999
$ws.Range("H125").Value = 37825
$ws.Range("I125").Value = 35650
$ws.Range("J125").Value = 40000
$ws.Range("K125").Value = 35650
$ws.Range("L125").Value = 40000
$ws.Range("M125").Value = -30730
$ws.Range("N125").Value = -49840
$ws.Range("H132").Value = 11064.125
$ws.Range("I132").Value = 11243.917
$ws.Range("K132").Value = 33731.751
$ws.Range("M132").Value = -31201.751
$ws.Range("H136").Value = 13170513
$ws.Range("I136").Value = 15632766
$ws.Range("J136").Value = 38499.168
$ws.Range("K136").Value = 46898298
$ws.Range("L136").Value = 115497.504
$ws.Range("M136").Value = -46895748
$ws.Range("N136").Value = -120597.504
